# Natmi following Dr Hou advice
# Extend the Rspo3-Sdc4 sending-cluster analysis to include "ECs" (not just
# "FAPs") as a sending cluster, producing 8 data rows (2 senders x 4 targets)
# instead of 4, and refresh the computed NATMI edge statistics accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Sdc4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06937033333333333
$ws.Range("H2").Value = 0.208111
$ws.Range("I2").Value = 0.01708561286819356
$ws.Range("J2").Value = 0.01708561286819356
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.001642333333333
$ws.Range("N2").Value = 9.004927
$ws.Range("O2").Value = 0.05169795991651582
$ws.Range("P2").Value = 0.05169795991651582
$ws.Range("Q2").Value = 0.2082249292107778
$ws.Range("R2").Value = 1.874024362897
$ws.Range("S2").Value = 0.0008832913292089773
$ws.Range("T2").Value = 0.0008832913292089776
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Sdc4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06937033333333333
$ws.Range("H3").Value = 0.208111
$ws.Range("I3").Value = 0.01708561286819356
$ws.Range("J3").Value = 0.01708561286819356
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.05428166666667
$ws.Range("N3").Value = 33.162845
$ws.Range("O3").Value = 0.1903903753498087
$ws.Range("P3").Value = 0.1903903753498088
$ws.Range("Q3").Value = 0.7668392039772223
$ws.Range("R3").Value = 6.901552835795001
$ws.Range("S3").Value = 0.003252936247056893
$ws.Range("T3").Value = 0.003252936247056895
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Sdc4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06937033333333333
$ws.Range("H4").Value = 0.208111
$ws.Range("I4").Value = 0.01708561286819356
$ws.Range("J4").Value = 0.01708561286819356
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 16.49405866666666
$ws.Range("N4").Value = 49.482176
$ws.Range("O4").Value = 0.2840808761059341
$ws.Range("P4").Value = 0.2840808761059341
$ws.Range("Q4").Value = 1.144198347726222
$ws.Range("R4").Value = 10.297785129536
$ws.Range("S4").Value = 0.004853695872403247
$ws.Range("T4").Value = 0.004853695872403248
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Sdc4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06937033333333333
$ws.Range("H5").Value = 0.208111
$ws.Range("I5").Value = 0.01708561286819356
$ws.Range("J5").Value = 0.01708561286819356
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.51115433333333
$ws.Range("N5").Value = 82.533463
$ws.Range("O5").Value = 0.4738307886277414
$ws.Range("P5").Value = 0.4738307886277414
$ws.Range("Q5").Value = 1.908457946488111
$ws.Range("R5").Value = 17.176121518393
$ws.Range("S5").Value = 0.00809568941952444
$ws.Range("T5").Value = 0.008095689419524441
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Sdc4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.990790333333333
$ws.Range("H6").Value = 11.972371
$ws.Range("I6").Value = 0.9829143871318063
$ws.Range("J6").Value = 0.9829143871318063
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.001642333333333
$ws.Range("N6").Value = 9.004927
$ws.Range("O6").Value = 0.05169795991651582
$ws.Range("P6").Value = 0.05169795991651582
$ws.Range("Q6").Value = 11.97892520799078
$ws.Range("R6").Value = 107.810326871917
$ws.Range("S6").Value = 0.05081466858730684
$ws.Range("T6").Value = 0.05081466858730684
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Sdc4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.990790333333333
$ws.Range("H7").Value = 11.972371
$ws.Range("I7").Value = 0.9829143871318063
$ws.Range("J7").Value = 0.9829143871318063
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.05428166666667
$ws.Range("N7").Value = 33.162845
$ws.Range("O7").Value = 0.1903903753498087
$ws.Range("P7").Value = 0.1903903753498088
$ws.Range("Q7").Value = 44.11532041727722
$ws.Range("R7").Value = 397.037883755495
$ws.Range("S7").Value = 0.1871374391027518
$ws.Range("T7").Value = 0.1871374391027519
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Rspo3"
$ws.Range("C8").Value = "Sdc4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.990790333333333
$ws.Range("H8").Value = 11.972371
$ws.Range("I8").Value = 0.9829143871318063
$ws.Range("J8").Value = 0.9829143871318063
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 16.49405866666666
$ws.Range("N8").Value = 49.482176
$ws.Range("O8").Value = 0.2840808761059341
$ws.Range("P8").Value = 0.2840808761059341
$ws.Range("Q8").Value = 65.82432988436621
$ws.Range("R8").Value = 592.4189689592959
$ws.Range("S8").Value = 0.2792271802335308
$ws.Range("T8").Value = 0.2792271802335309
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Rspo3"
$ws.Range("C9").Value = "Sdc4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.990790333333333
$ws.Range("H9").Value = 11.972371
$ws.Range("I9").Value = 0.9829143871318063
$ws.Range("J9").Value = 0.9829143871318063
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.51115433333333
$ws.Range("N9").Value = 82.533463
$ws.Range("O9").Value = 0.4738307886277414
$ws.Range("P9").Value = 0.4738307886277414
$ws.Range("Q9").Value = 109.7912487723081
$ws.Range("R9").Value = 988.1212389507729
$ws.Range("S9").Value = 0.4657350992082169
$ws.Range("T9").Value = 0.4657350992082169
